$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old pin code in A1 and add three more pin codes below it.
# Leading apostrophe forces these numeric-looking values to be stored as
# text (shared strings), matching the original "744105" text cell.
$ws.Range("A1").Value = "'249131"
$ws.Range("A2").Value = "'249171"
$ws.Range("A3").Value = "'12404"
$ws.Range("A4").Value = "'12405"

# The apostrophe prefix otherwise leaves a "quote prefix" number format on
# the cells; clear formatting so the cells keep using the default style
# (same as the original A1 cell, style index 0).
$ws.Range("A1:A4").ClearFormats()
